$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.522.96"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").Value = "1.728.58"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.94"
$ws.Range("E5").Value = "  +1.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4808"
$ws.Range("E7").Value = "  +1.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2675"
$ws.Range("E8").Value = "  +1.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06188"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").Value = "1.728.70"
$ws.Range("E10").Value = "  +0.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07195"
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.58"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6105"
$ws.Range("E13").Value = "  +2.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.531"
$ws.Range("E14").Value = "  +2.26%  "
$ws.Range("E15").Value = "  +1.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9998"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "26.517.88"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9998"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006943"
$ws.Range("E19").Value = "  +1.75%  "
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "1.952.08"
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.529"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.811"
$ws.Range("E23").Value = "  +0.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.257"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.90"
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("E26").Value = "  +0.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.780"
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.406"
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "107.32"
$ws.Range("E29").Value = "  +0.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.983"
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08031"
$ws.Range("E31").Value = "  +2.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.697"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04517"
$ws.Range("E33").Value = "  +0.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.618"
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("E35").Value = "  +1.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6256"
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9104"
$ws.Range("E37").Value = "  -1.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.072"
$ws.Range("E38").Value = "  +7.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.389"
$ws.Range("E39").Value = "  -2.59%  "
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.63"
$ws.Range("E41").Value = "  -9.63%  "
$ws.Range("E42").Value = "  +1.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.537"
$ws.Range("E43").Value = "  -2.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3881"
$ws.Range("E44").Value = "  +1.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.991"
$ws.Range("E45").Value = "  +9.98%  "
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05370"
$ws.Range("E47").Value = "  +1.91%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.65"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.807"
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.250"
$ws.Range("E50").Value = "  +2.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3401"
$ws.Range("E51").Value = "  +0.40%  "
